# Update the cryptos list worksheet with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)

# Row 2 - Bitcoin
$ws.Range("D2").Value = "53.947.78"
$ws.Range("E2").Value = "  -2.75%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.284.21"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.42%  "

# Row 5 - BNB
$ws.Range("D5").Value = "493.89"
$ws.Range("E5").Value = "  -1.09%  "

# Row 6 - Solana
$ws.Range("D6").Value = "127.89"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.62%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.32%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.281.87"
$ws.Range("E9").Value = "  -1.61%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0934"
$ws.Range("E10").Value = "  -3.95%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.44%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.06%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "4.67"
$ws.Range("E13").Value = "  -3.22%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.686.72"
$ws.Range("E14").Value = "  -2.71%  "

# Row 15 - Avalanche
$ws.Range("E15").Value = "  +0.46%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "53.981.87"
$ws.Range("E16").Value = "  -3.03%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -1.69%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.264.06"
$ws.Range("E18").Value = "  -3.23%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "9.95"
$ws.Range("E19").Value = "  +1.01%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "4.04"
$ws.Range("E20").Value = "  +1.89%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "298.97"
$ws.Range("E21").Value = "  -2.47%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +1.50%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.01%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "63.51"
$ws.Range("E24").Value = "  -2.39%  "

# Row 25 - Binance-PegBSC-USD
$ws.Range("E25").Value = "  +0.02%  "

# Row 26 - Polygon
$ws.Range("E26").Value = "  +1.05%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "2.383.76"
$ws.Range("E27").Value = "  -3.79%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.73%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  +0.99%  "

# Row 30 - Monero
$ws.Range("D30").Value = "162.94"
$ws.Range("E30").Value = "  -5.66%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.46%  "

# Row 32 - PEPE
$ws.Range("D32").Value = "0.0₃0684"
$ws.Range("E32").Value = "  -1.60%  "

# Row 33/34 - swap USDe <-> Aptos
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "5.83"
$ws.Range("E33").Value = "  +1.38%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.22%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  -0.67%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  +0.80%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "'17.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.29%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +1.66%  "

# Row 39 - SuiNetwork
$ws.Range("D39").Value = "0.859"
$ws.Range("E39").Value = "  +4.97%  "

# Row 40 - NEARProtocol
$ws.Range("D40").Value = "3.65"
$ws.Range("E40").Value = "  +1.06%  "

# Row 41 - OKB
$ws.Range("D41").Value = "35.32"
$ws.Range("E41").Value = "  -2.09%  "

# Row 42 - PolygonEcosystemToken
$ws.Range("D42").Value = "0.375"
$ws.Range("E42").Value = "  +1.34%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +1.92%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +0.25%  "

# Row 45/46 - swap Aave <-> RenderToken
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "4.88"
$ws.Range("E45").Value = "  +4.80%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "126.89"
$ws.Range("E46").Value = "  +0.68%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  -0.05%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  -1.86%  "

# Row 49 - Bittensor
$ws.Range("D49").Value = "238.95"
$ws.Range("E49").Value = "  +2.13%  "

# Row 50 - Hedera
$ws.Range("E50").Value = "  +1.11%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -0.27%  "
